$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.803.18"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.889.67"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "365.15"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("D6").Value = "101.17"
$ws.Range("E6").Value = "  -4.39%  "
$ws.Range("D7").Value = "0.539"
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.577"
$ws.Range("E9").Value = "  -4.57%  "
$ws.Range("D10").Value = "36.02"
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "0.0823"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.344.63"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "18.05"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").Value = "7.27"
$ws.Range("E15").Value = "  -3.99%  "
$ws.Range("D16").Value = "2.891.88"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "0.910"
$ws.Range("E17").Value = "  -5.00%  "
$ws.Range("D18").Value = "50.794.72"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "3.18"
$ws.Range("E19").Value = "  -6.27%  "
$ws.Range("D20").Value = "7.09"
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").Value = "12.64"
$ws.Range("E21").Value = "  -5.54%  "
$ws.Range("D22").Value = "0.0₃0933"
$ws.Range("E22").Value = "  -3.16%  "
$ws.Range("D23").Value = "67.57"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").Value = "255.88"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "2.64"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").Value = "4.32"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("D29").Value = "25.28"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").Value = "6.85"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("D31").Value = "0.100"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("D32").Value = "6.08"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").Value = "9.76"
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("E34").Value = "  -1.41%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "50.78"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "33.68"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").Value = "0.0415"
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("D40").Value = "2.58"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "16.74"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  -6.03%  "
$ws.Range("D43").Value = "0.111"
$ws.Range("E43").Value = "  -3.73%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "118.18"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "21.60"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.001.40"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").Value = "3.10"
$ws.Range("E49").Value = "  -5.94%  "
$ws.Range("D50").Value = "3.184.52"
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("D51").Value = "0.233"
$ws.Range("E51").Value = "  -2.25%  "
